$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 6 new blank rows before row 34 (the former "last row" of the
#    worker table), shifting the existing rows 34-39 (footer/signature rows)
#    down to 40-45. This grows the 3-worker x 6-period table (18 rows,
#    16..33) into a 4-worker x 6-period table (24 rows, 16..39).
# ---------------------------------------------------------------------------
$ws.Rows("34:39").Insert(-4121) | Out-Null   # xlShiftDown

# Give the newly inserted rows (34-38) the same formatting as the existing
# "normal" data rows (copy format only, from row 32).
$ws.Range("B32:J32").Copy() | Out-Null
$ws.Range("B34:J38").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# The new last row of the table (39) should carry the special bottom-border
# formatting that used to belong to row 33.
$ws.Range("B33:J33").Copy() | Out-Null
$ws.Range("B39:J39").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 33 itself is no longer the last row of the table, so it now needs the
# regular "normal" row formatting instead of the special bottom border.
$ws.Range("B32:J32").Copy() | Out-Null
$ws.Range("B33:J33").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Rewrite the whole worker table (rows 16-39) with the refreshed data:
#    previous statements removed, new ones added (per commit message),
#    grouped by worker, periods descending 2412 -> 2407.
# ---------------------------------------------------------------------------
$workers = @(
    @{ Doc = "3809494";     Nombre = "JAIRO ENRIQUE MUENTES ALEAN" },
    @{ Doc = "73353711";    Nombre = "JUAN ANTONIO MATA ACEVEDO" },
    @{ Doc = "1143354628";  Nombre = "DAVID JOSE ESCOBAR LOPEZ" },
    @{ Doc = "1047504382";  Nombre = "JESWALDO MONTERROZA GONZALEZ" }
)
$periodos = @("2412", "2411", "2410", "2409", "2408", "2407")

$r = 16
foreach ($w in $workers) {
    for ($i = 0; $i -lt $periodos.Length; $i++) {
        $valor = 52000
        if ($i -eq 0) { $valor = 19067 }

        $ws.Range("B$r").Value2 = "CC"
        $ws.Range("C$r").Value2 = $w.Doc
        $ws.Range("D$r").Value2 = $w.Nombre
        $ws.Range("E$r").Value2 = $periodos[$i]
        $ws.Range("F$r").Value2 = $valor
        $ws.Range("G$r").Value2 = 1300000

        $r = $r + 1
    }
}

# ---------------------------------------------------------------------------
# 3. Update the summary header cells.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 1116268   # VALOR MORA (total)
$ws.Range("C13").Value2 = 4         # Cant. Trabajadores
